$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r1 = New-Object 'object[,]' 1,4
$r1[0,0] = 5.270858666666666
$r1[0,1] = 15.812576
$r1[0,2] = 0.02659672411376159
$r1[0,3] = 0.02659672411376159
$ws.Range("G2:J2").Value = $r1
$r2 = New-Object 'object[,]' 1,8
$r2[0,0] = 0.353079
$r2[0,1] = 1.059237
$r2[0,2] = 0.01390801122570493
$r2[0,3] = 0.01390801122570493
$r2[0,4] = 1.861029507168
$r2[0,5] = 16.749265564512
$r2[0,6] = 0.0003699075375411731
$r2[0,7] = 0.0003699075375411732
$ws.Range("M2:T2").Value = $r2

$r3 = New-Object 'object[,]' 1,4
$r3[0,0] = 5.270858666666666
$r3[0,1] = 15.812576
$r3[0,2] = 0.02659672411376159
$r3[0,3] = 0.02659672411376159
$ws.Range("G3:J3").Value = $r3
$r4 = New-Object 'object[,]' 1,8
$r4[0,0] = 24.359095
$r4[0,1] = 73.077285
$r4[0,2] = 0.9595205795530543
$r4[0,3] = 0.9595205795530543
$r4[0,4] = 128.3933469929067
$r4[0,5] = 1155.54012293616
$r4[0,6] = 0.02552010413584921
$r4[0,7] = 0.02552010413584921
$ws.Range("M3:T3").Value = $r4

$r5 = New-Object 'object[,]' 1,14
$r5[0,0] = 5.270858666666666
$r5[0,1] = 15.812576
$r5[0,2] = 0.02659672411376159
$r5[0,3] = 0.02659672411376159
$r5[0,4] = 3
$r5[0,5] = 1
$r5[0,6] = 0.6745613333333332
$r5[0,7] = 2.023684
$r5[0,8] = 0.02657140922124081
$r5[0,9] = 0.02657140922124081
$r5[0,10] = 3.555517449998221
$r5[0,11] = 31.999657049984
$r5[0,12] = 0.0007067124403712026
$r5[0,13] = 0.0007067124403712026
$ws.Range("G4:T4").Value = $r5

$r6 = New-Object 'object[,]' 1,2
$r6[0,0] = 0.01269587703542364
$r6[0,1] = 0.01269587703542364
$ws.Range("I5:J5").Value = $r6
$r7 = New-Object 'object[,]' 1,8
$r7[0,0] = 0.353079
$r7[0,1] = 1.059237
$r7[0,2] = 0.01390801122570493
$r7[0,3] = 0.01390801122570493
$r7[0,4] = 0.8883575917559999
$r7[0,5] = 7.995218325803999
$r7[0,6] = 0.0001765744003288413
$r7[0,7] = 0.0001765744003288413
$ws.Range("M5:T5").Value = $r7

$r8 = New-Object 'object[,]' 1,2
$r8[0,0] = 0.01269587703542364
$r8[0,1] = 0.01269587703542364
$ws.Range("I6:J6").Value = $r8
$r9 = New-Object 'object[,]' 1,8
$r9[0,0] = 24.359095
$r9[0,1] = 73.077285
$r9[0,2] = 0.9595205795530543
$r9[0,3] = 0.9595205795530543
$r9[0,4] = 61.28823003224667
$r9[0,5] = 551.59407029022
$r9[0,6] = 0.012181955290964
$r9[0,7] = 0.012181955290964
$ws.Range("M6:T6").Value = $r9

$r10 = New-Object 'object[,]' 1,12
$r10[0,0] = 0.01269587703542364
$r10[0,1] = 0.01269587703542364
$r10[0,2] = 3
$r10[0,3] = 1
$r10[0,4] = 0.6745613333333332
$r10[0,5] = 2.023684
$r10[0,6] = 0.02657140922124081
$r10[0,7] = 0.02657140922124081
$r10[0,8] = 1.697217001214222
$r10[0,9] = 15.274953010928
$r10[0,10] = 0.000337347344130795
$r10[0,11] = 0.000337347344130795
$ws.Range("I7:T7").Value = $r10

$r11 = New-Object 'object[,]' 1,4
$r11[0,0] = 114.018682
$r11[0,1] = 342.056046
$r11[0,2] = 0.5753376481419691
$r11[0,3] = 0.5753376481419691
$ws.Range("G8:J8").Value = $r11
$r12 = New-Object 'object[,]' 1,8
$r12[0,0] = 0.353079
$r12[0,1] = 1.059237
$r12[0,2] = 0.01390801122570493
$r12[0,3] = 0.01390801122570493
$r12[0,4] = 40.257602221878
$r12[0,5] = 362.318419996902
$r12[0,6] = 0.008001802468929179
$r12[0,7] = 0.008001802468929179
$ws.Range("M8:T8").Value = $r12

$r13 = New-Object 'object[,]' 1,4
$r13[0,0] = 114.018682
$r13[0,1] = 342.056046
$r13[0,2] = 0.5753376481419691
$r13[0,3] = 0.5753376481419691
$ws.Range("G9:J9").Value = $r13
$r14 = New-Object 'object[,]' 1,8
$r14[0,0] = 24.359095
$r14[0,1] = 73.077285
$r14[0,2] = 0.9595205795530543
$r14[0,3] = 0.9595205795530543
$r14[0,4] = 2777.39190661279
$r14[0,5] = 24996.52715951511
$r14[0,6] = 0.5520483135838734
$r14[0,7] = 0.5520483135838734
$ws.Range("M9:T9").Value = $r14

$r15 = New-Object 'object[,]' 1,14
$r15[0,0] = 114.018682
$r15[0,1] = 342.056046
$r15[0,2] = 0.5753376481419691
$r15[0,3] = 0.5753376481419691
$r15[0,4] = 3
$r15[0,5] = 1
$r15[0,6] = 0.6745613333333332
$r15[0,7] = 2.023684
$r15[0,8] = 0.02657140922124081
$r15[0,9] = 0.02657140922124081
$r15[0,10] = 76.91259415482934
$r15[0,11] = 692.213347393464
$r15[0,12] = 0.01528753208916652
$r15[0,13] = 0.01528753208916652
$ws.Range("G10:T10").Value = $r15

$r16 = New-Object 'object[,]' 1,4
$r16[0,0] = 1.265015666666667
$r16[0,1] = 3.795047
$r16[0,2] = 0.006383262161570549
$r16[0,3] = 0.006383262161570549
$ws.Range("G11:J11").Value = $r16
$r17 = New-Object 'object[,]' 1,8
$r17[0,0] = 0.353079
$r17[0,1] = 1.059237
$r17[0,2] = 0.01390801122570493
$r17[0,3] = 0.01390801122570493
$r17[0,4] = 0.446650466571
$r17[0,5] = 4.019854199139
$r17[0,6] = 0.00008877848179974071
$r17[0,7] = 0.00008877848179974071
$ws.Range("M11:T11").Value = $r17

$r18 = New-Object 'object[,]' 1,4
$r18[0,0] = 1.265015666666667
$r18[0,1] = 3.795047
$r18[0,2] = 0.006383262161570549
$r18[0,3] = 0.006383262161570549
$ws.Range("G12:J12").Value = $r18
$r19 = New-Object 'object[,]' 1,8
$r19[0,0] = 24.359095
$r19[0,1] = 73.077285
$r19[0,2] = 0.9595205795530543
$r19[0,3] = 0.9595205795530543
$r19[0,4] = 30.81463680082167
$r19[0,5] = 277.331731207395
$r19[0,6] = 0.006124871408709255
$r19[0,7] = 0.006124871408709255
$ws.Range("M12:T12").Value = $r19

$r20 = New-Object 'object[,]' 1,14
$r20[0,0] = 1.265015666666667
$r20[0,1] = 3.795047
$r20[0,2] = 0.006383262161570549
$r20[0,3] = 0.006383262161570549
$r20[0,4] = 3
$r20[0,5] = 1
$r20[0,6] = 0.6745613333333332
$r20[0,7] = 2.023684
$r20[0,8] = 0.02657140922124081
$r20[0,9] = 0.02657140922124081
$r20[0,10] = 0.8533306547942222
$r20[0,11] = 7.679975893148
$r20[0,12] = 0.0001696122710615533
$r20[0,13] = 0.0001696122710615533
$ws.Range("G13:T13").Value = $r20

$r21 = New-Object 'object[,]' 1,4
$r21[0,0] = 75.10640066666666
$r21[0,1] = 225.319202
$r21[0,2] = 0.3789864885472752
$r21[0,3] = 0.3789864885472752
$ws.Range("G14:J14").Value = $r21
$r22 = New-Object 'object[,]' 1,8
$r22[0,0] = 0.353079
$r22[0,1] = 1.059237
$r22[0,2] = 0.01390801122570493
$r22[0,3] = 0.01390801122570493
$r22[0,4] = 26.518492840986
$r22[0,5] = 238.666435568874
$r22[0,6] = 0.005270948337105995
$r22[0,7] = 0.005270948337105996
$ws.Range("M14:T14").Value = $r22

$r23 = New-Object 'object[,]' 1,4
$r23[0,0] = 75.10640066666666
$r23[0,1] = 225.319202
$r23[0,2] = 0.3789864885472752
$r23[0,3] = 0.3789864885472752
$ws.Range("G15:J15").Value = $r23
$r24 = New-Object 'object[,]' 1,8
$r24[0,0] = 24.359095
$r24[0,1] = 73.077285
$r24[0,2] = 0.9595205795530543
$r24[0,3] = 0.9595205795530543
$r24[0,4] = 1829.523948947397
$r24[0,5] = 16465.71554052657
$r24[0,6] = 0.3636453351336584
$r24[0,7] = 0.3636453351336584
$ws.Range("M15:T15").Value = $r24

$r25 = New-Object 'object[,]' 1,14
$r25[0,0] = 75.10640066666666
$r25[0,1] = 225.319202
$r25[0,2] = 0.3789864885472752
$r25[0,3] = 0.3789864885472752
$r25[0,4] = 3
$r25[0,5] = 1
$r25[0,6] = 0.6745613333333332
$r25[0,7] = 2.023684
$r25[0,8] = 0.02657140922124081
$r25[0,9] = 0.02657140922124081
$r25[0,10] = 50.66387377557421
$r25[0,11] = 455.9748639801679
$r25[0,12] = 0.01007020507651074
$r25[0,13] = 0.01007020507651074
$ws.Range("G16:T16").Value = $r25
